# 9.3.2.xlsx - add the 2023 data point (column S) to the existing
# 2008-2022 (columns D-R) time series on the single worksheet.
#
# Row 4 holds the year headers, row 5 the "small-scale industries with a
# loan" series, row 6 the "Industry" series. Column R already holds 2022 /
# the matching 2021-2022 figures, so column S is the next empty column and
# simply continues each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing R-column formatting (font/border/number format) for
# each row onto the new S cells, then overwrite with the new values - this
# keeps S4/S5/S6 visually identical to R4/R5/R6 instead of falling back to
# the sheet's default (unstyled) look.
$ws.Range("R4").Copy($ws.Range("S4"))
$ws.Range("S4").Value = 2023

$ws.Range("R5").Copy($ws.Range("S5"))
$ws.Range("S5").Value = 7.1262361838278068

$ws.Range("R6").Copy($ws.Range("S6"))
$ws.Range("S6").Value = 10.974456007568591

# Reset the active selection back to the top-left corner now that S4 holds
# real data (it was previously parked on S4 awaiting this value).
$ws.Range("A1").Select()
